# Update cryptos list values per latest scrape (Coin/Link/Price/Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.670.76'
$ws.Range("E2").Value = '''  -1.41%  '
$ws.Range("D3").Value = '''1.794.13'
$ws.Range("E3").Value = '''  -1.33%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("D5").Value = '''309.27'
$ws.Range("E5").Value = '''  -0.41%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '''  +0.00%  '
$ws.Range("D7").Value = '''0.4458'
$ws.Range("E7").Value = '''  +5.75%  '
$ws.Range("D8").Value = '''0.3678'
$ws.Range("E8").Value = '''  +0.28%  '
$ws.Range("D9").Value = '''0.07330'
$ws.Range("E9").Value = '''  +1.67%  '
$ws.Range("D10").Value = '''0.8578'
$ws.Range("E10").Value = '''  +1.78%  '
$ws.Range("D11").Value = '''20.60'
$ws.Range("E11").Value = '''  -1.08%  '
$ws.Range("D12").Value = '''1.796.15'
$ws.Range("E12").Value = '''  -1.40%  '
$ws.Range("D13").Value = '''6.620'
$ws.Range("E13").Value = '''  -0.21%  '
$ws.Range("D14").Value = '''92.03'
$ws.Range("E14").Value = '''  +3.48%  '
$ws.Range("D15").Value = '''0.07070'
$ws.Range("E15").Value = '''  -0.07%  '
$ws.Range("D16").Value = '''5.261'
$ws.Range("E16").Value = '''  -0.13%  '
$ws.Range("D17").Value = '''1.003'
$ws.Range("E17").Value = '''  -0.04%  '
$ws.Range("D18").Value = '''0.000008667'
$ws.Range("E18").Value = '''  -1.64%  '
$ws.Range("D19").Value = '''1.001'
$ws.Range("E19").Value = '''  +0.03%  '
$ws.Range("D20").Value = '''14.79'
$ws.Range("E20").Value = '''  -1.16%  '
$ws.Range("D21").Value = '''26.700.09'
$ws.Range("E21").Value = '''  -1.51%  '
$ws.Range("D22").Value = '''5.155'
$ws.Range("E22").Value = '''  +0.82%  '
$ws.Range("D23").Value = '''10.78'
$ws.Range("E23").Value = '''  -0.36%  '
$ws.Range("D24").Value = '''1.980'
$ws.Range("E24").Value = '''  +0.08%  '
$ws.Range("D25").Value = '''151.88'
$ws.Range("D26").Value = '''18.42'
$ws.Range("E26").Value = '''  +0.85%  '
$ws.Range("D27").Value = '''2.169'
$ws.Range("E27").Value = '''  -3.20%  '
$ws.Range("D28").Value = '''5.175'
$ws.Range("E28").Value = '''  -0.41%  '
$ws.Range("D29").Value = '''117.29'
$ws.Range("E29").Value = '''  +1.12%  '
$ws.Range("D30").Value = '''0.08773'
$ws.Range("E30").Value = '''  -0.26%  '
$ws.Range("D31").Value = '''0.7385'
$ws.Range("E31").Value = '''  +0.02%  '
$ws.Range("D32").Value = '''1.153'
$ws.Range("E32").Value = '''  -2.06%  '
$ws.Range("D33").Value = '''2.918'
$ws.Range("E33").Value = '''  -1.90%  '
$ws.Range("D34").Value = '''4.443'
$ws.Range("E34").Value = '''  +0.74%  '
$ws.Range("D35").Value = '''1.001'
$ws.Range("E35").Value = '''  -0.07%  '
$ws.Range("E36").Value = '''  -1.09%  '
$ws.Range("D37").Value = '''0.01955'
$ws.Range("E37").Value = '''  -0.40%  '
$ws.Range("D38").Value = '''0.05182'
$ws.Range("E38").Value = '''  -1.24%  '
$ws.Range("D39").Value = '''0.5256'
$ws.Range("E39").Value = '''  +4.64%  '
$ws.Range("E40").Value = '''  -1.38%  '
$ws.Range("D41").Value = '''6.955'
$ws.Range("E41").Value = '''  -4.16%  '
$ws.Range("E42").Value = '''  -0.38%  '
$ws.Range("D43").Value = '''0.5067'
$ws.Range("E43").Value = '''  +6.71%  '
$ws.Range("D44").Value = '''8.404'
$ws.Range("E44").Value = '''  -1.71%  '
$ws.Range("D45").Value = '''1.963'
$ws.Range("E45").Value = '''  +4.57%  '
$ws.Range("E46").Value = '''  -1.06%  '
$ws.Range("D47").Value = '''105.00'
$ws.Range("E47").Value = '''  -0.94%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '''1.001'
$ws.Range("E48").Value = '''  -0.02%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.667'
$ws.Range("E49").Value = '''  +1.05%  '
$ws.Range("D50").Value = '''0.06295'
$ws.Range("E50").Value = '''  -1.19%  '
$ws.Range("D51").Value = '''0.9140'
$ws.Range("E51").Value = '''  +1.18%  '
